# TC42_Canine_Filter_Breed-ShihTzu.xlsx — "startup" sheet, row 4 (FilesTab)
# column B holds the Neo4j query for the Case Files tab. The query's RETURN
# clause column order is being fixed (File Name, Format, File Type, Size,
# Association, ... instead of File Name, File Type, Association,
# Description, Format, Size, ...).
#
# The replacement text is base64-encoded to avoid any PowerShell quoting /
# escaping issues with the backticks, single quotes and newlines it
# contains.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$b64 = "TUFUQ0ggKGY6ZmlsZSktLT4ocGFyZW50KQpXSVRIIERJU1RJTkNUIGYsIHBhcmVudApNQVRDSCAoZGlhZzpkaWFnbm9zaXMpLS0+KGMpCk9QVElPTkFMIE1BVENIIChmKS1bKl0tPihzYW1wOnNhbXBsZSkKTUFUQ0ggKGYpLVsqXS0+KGM6Y2FzZSk8LS0oZGVtbzpkZW1vZ3JhcGhpYykKV0hFUkUgZGVtby5icmVlZCBJTiAgWydTaGloIFR6dSddIApPUFRJT05BTCBNQVRDSCAoczpzdHVkeSk8LS0oYyk8LS0oZGlhZzpkaWFnbm9zaXMpPC1bKl0tKHNhbXApCldJVEgKICAgICAgICBmLCBwYXJlbnQsIGMsIGRlbW8sIGRpYWcsIHMsIHNhbXAsCiAgICAgICAgWydCeXRlcycsICdLQicsICdNQicsICdHQicsICdUQiddIEFTIHVuaXRzLAogICAgICAgIHRvSW50ZWdlcihmbG9vcihsb2coZi5maWxlX3NpemUpL2xvZygxMDI0KSkpIGFzIGksCiAgICAgICAgMiBhcyBwcmVjaXNpb24KV0lUSAogICAgICAgIGYsIHBhcmVudCwgYywgZGVtbywgZGlhZywgcywgc2FtcCwKICAgICAgICBmLmZpbGVfc2l6ZSAvKDEwMjReaSkgQVMgdmFsdWUsIAogICAgICAgIDEwXnByZWNpc2lvbiBBUyBmYWN0b3IsCiAgICAgICAgdW5pdHNbaV0gYXMgdW5pdApXSVRIICAgIAogICAgICAgIGYsIHBhcmVudCwgYywgZGVtbywgZGlhZywgcywgc2FtcCwgdW5pdCwKICAgICAgICByb3VuZChmYWN0b3IgKiB2YWx1ZSkvZmFjdG9yIEFTIHNpemUKUkVUVVJOIAogICAgICAgY29hbGVzY2UoZi5maWxlX25hbWUsICcnKSBBUyBgRmlsZSBOYW1lYCwKICAgICAgIGNvYWxlc2NlKGYuZmlsZV9mb3JtYXQsICcnKSBBUyBgRm9ybWF0YCwKICAgICAgICBjb2FsZXNjZShmLmZpbGVfdHlwZSwgJycpIEFTIGBGaWxlIFR5cGVgLAogICAgICAgQ0FTRSBzaXplICUgMSBXSEVOIDAgVEhFTiBhcG9jLmNvbnZlcnQudG9JbnRlZ2VyKHNpemUpKycgJyArdW5pdCBFTFNFIHNpemUrJyAnICt1bml0IEVORCBBUyBTaXplLAogICAgICAgIGNvYWxlc2NlKGxhYmVscyhwYXJlbnQpWzBdLCAnJykgQVMgYEFzc29jaWF0aW9uYCwKICAgICAgICBjb2FsZXNjZShmLmZpbGVfZGVzY3JpcHRpb24sICcnKSBBUyBgRGVzY3JpcHRpb25gLAogICAgICAgIGNvYWxlc2NlKHNhbXAuc2FtcGxlX2lkLCAnJykgQVMgYFNhbXBsZSBJRGAsCiAgICAgICAgY29hbGVzY2UoYy5jYXNlX2lkLCAnJykgQVMgYENhc2UgSURgLAogICAgICAgIGNvYWxlc2NlKGRlbW8uYnJlZWQsJycpIEFTIEJyZWVkICwKICAgICAgICBjb2FsZXNjZShkaWFnLmRpc2Vhc2VfdGVybSwnJykgQVMgRGlhZ25vc2lz"
$bytes = [System.Convert]::FromBase64String($b64)
$newQuery = [System.Text.Encoding]::UTF8.GetString($bytes)

$ws.Range("B4").Value2 = $newQuery

Write-Output $ws.Range("B4").Value2
